$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.067.39'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '1.818.58'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("D5").Value = '''232.76'
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").Value = '''0.5862'
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("D8").Value = '''0.2722'
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("D9").Value = '''0.06768'
$ws.Range("E9").Value = '  -4.41%  '
$ws.Range("D10").Value = '''22.90'
$ws.Range("D11").Value = '''0.07514'
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").Value = '1.812.67'
$ws.Range("E12").Value = '  -0.96%  '
$ws.Range("D13").Value = '''4.637'
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '''0.6201'
$ws.Range("E14").Value = '  -2.77%  '
$ws.Range("D15").Value = '''0.000009165'
$ws.Range("E15").Value = '  -7.94%  '
$ws.Range("D16").Value = '''74.35'
$ws.Range("E16").Value = '  -6.84%  '
$ws.Range("D17").Value = '28.805.64'
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("D18").Value = '''5.393'
$ws.Range("E18").Value = '  -9.98%  '
$ws.Range("D20").Value = '''206.76'
$ws.Range("E20").Value = '  -10.36%  '
$ws.Range("D21").Value = '''11.32'
$ws.Range("E21").Value = '  -4.15%  '
$ws.Range("D22").Value = '''6.736'
$ws.Range("E22").Value = '  -4.16%  '
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '''153.90'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").Value = '''7.742'
$ws.Range("E25").Value = '  -4.18%  '
$ws.Range("D26").Value = '''0.1257'
$ws.Range("E26").Value = '  -2.84%  '
$ws.Range("D27").Value = '''16.17'
$ws.Range("E27").Value = '  -3.49%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''1.410'
$ws.Range("E28").Value = '  -3.96%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '''0.06295'
$ws.Range("E29").Value = '  -6.26%  '
$ws.Range("D30").Value = '''1.426'
$ws.Range("E30").Value = '  -2.43%  '
$ws.Range("D31").Value = '''3.678'
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("D32").Value = '''3.648'
$ws.Range("E32").Value = '  -5.12%  '
$ws.Range("D33").Value = '''1.677'
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("D34").Value = '''1.040'
$ws.Range("E34").Value = '  -8.07%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.536'
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''0.6323'
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("D37").Value = '''2.744'
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = '''0.01694'
$ws.Range("E38").Value = '  -4.17%  '
$ws.Range("D39").Value = '''6.384'
$ws.Range("D40").Value = '1.124.90'
$ws.Range("E40").Value = '  -9.15%  '
$ws.Range("D41").Value = '''0.8610'
$ws.Range("E41").Value = '  -7.24%  '
$ws.Range("D42").Value = '''1.006'
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("D43").Value = '1.967.62'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '''99.60'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").Value = '''59.81'
$ws.Range("E45").Value = '  -6.15%  '
$ws.Range("E46").Value = '  -3.15%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05483'
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '''0.4528'
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.561'
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("D50").Value = '''1.010'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("D51").Value = '''8.150'
$ws.Range("E51").Value = '  -4.50%  '
